$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row: new columns I and J
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Data rows 2-27: I = 1, J = same as H
for ($r = 2; $r -le 27; $r++) {
    $hVal = $ws.Cells.Item($r, 8).Value2
    $ws.Cells.Item($r, 9).Value = 1
    $ws.Cells.Item($r, 10).Value = $hVal
}
